$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.282.92'
$ws.Range("E2").Value = '  -4.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.538.88'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.17'
$ws.Range("E5").Value = '  -4.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.17'
$ws.Range("E6").Value = '  -8.05%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -4.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.542.28'
$ws.Range("E9").Value = '  -4.44%  '
$ws.Range("E10").Value = '  -8.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("E12").Value = '  -5.97%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.987.21'
$ws.Range("E14").Value = '  -4.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.265.59'
$ws.Range("E15").Value = '  -4.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.64'
$ws.Range("E16").Value = '  -6.53%  '
$ws.Range("E17").Value = '  -6.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.548.44'
$ws.Range("E18").Value = '  -4.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.50'
$ws.Range("E19").Value = '  -5.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.79'
$ws.Range("E20").Value = '  -3.96%  '
$ws.Range("E21").Value = '  -6.05%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.92'
$ws.Range("E23").Value = '  -5.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.49'
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.657.45'
$ws.Range("E27").Value = '  -3.70%  '
$ws.Range("E28").Value = '  -6.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0781'
$ws.Range("E29").Value = '  -9.37%  '
$ws.Range("E30").Value = '  -6.44%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '148.76'
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.45'
$ws.Range("E33").Value = '  -5.63%  '
$ws.Range("E34").Value = '  -6.27%  '
$ws.Range("E36").Value = '  -7.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.898'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  -8.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.87'
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("E40").Value = '  -11.26%  '
$ws.Range("E41").Value = '  -8.14%  '
$ws.Range("E42").Value = '  -8.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '279.47'
$ws.Range("E43").Value = '  -9.57%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0989'
$ws.Range("E45").Value = '  -3.10%  '
$ws.Range("E46").Value = '  -7.58%  '
$ws.Range("E47").Value = '  -5.98%  '
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.55'
$ws.Range("E49").Value = '  -6.90%  '
$ws.Range("E50").Value = '  -6.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.50'
$ws.Range("E51").Value = '  -11.00%  '
